# Replace each arithmetic answer/equation string in the table cells with
# the new equation (commit: "Update master to output generated at c8c62b6").
$d = $word.ActiveDocument

$d.Content.Find.Execute("11+69=80", $true, $true, $false, $false, $false, $true, 1, $false, "15+54=69", 2) | Out-Null
$d.Content.Find.Execute("86-15=71", $true, $true, $false, $false, $false, $true, 1, $false, "29-6=23", 2) | Out-Null
$d.Content.Find.Execute("75-44=31", $true, $true, $false, $false, $false, $true, 1, $false, "65-55=10", 2) | Out-Null
$d.Content.Find.Execute("79-1=78", $true, $true, $false, $false, $false, $true, 1, $false, "51-7=44", 2) | Out-Null
$d.Content.Find.Execute("26+34=60", $true, $true, $false, $false, $false, $true, 1, $false, "93+2=95", 2) | Out-Null
$d.Content.Find.Execute("96-57=39", $true, $true, $false, $false, $false, $true, 1, $false, "76-30=46", 2) | Out-Null
$d.Content.Find.Execute("41+37=78", $true, $true, $false, $false, $false, $true, 1, $false, "77-21=56", 2) | Out-Null
$d.Content.Find.Execute("35+11=46", $true, $true, $false, $false, $false, $true, 1, $false, "65-47=18", 2) | Out-Null
$d.Content.Find.Execute("51-50=1", $true, $true, $false, $false, $false, $true, 1, $false, "29-10=19", 2) | Out-Null
$d.Content.Find.Execute("86+10=96", $true, $true, $false, $false, $false, $true, 1, $false, "18+44=62", 2) | Out-Null
$d.Content.Find.Execute("39-4=35", $true, $true, $false, $false, $false, $true, 1, $false, "7+10=17", 2) | Out-Null
$d.Content.Find.Execute("21-10=11", $true, $true, $false, $false, $false, $true, 1, $false, "42-18=24", 2) | Out-Null
$d.Content.Find.Execute("78-34=44", $true, $true, $false, $false, $false, $true, 1, $false, "66-56=10", 2) | Out-Null
$d.Content.Find.Execute("5+73=78", $true, $true, $false, $false, $false, $true, 1, $false, "87-68=19", 2) | Out-Null
$d.Content.Find.Execute("19+38=57", $true, $true, $false, $false, $false, $true, 1, $false, "13+57=70", 2) | Out-Null
$d.Content.Find.Execute("42+2=44", $true, $true, $false, $false, $false, $true, 1, $false, "25+55=80", 2) | Out-Null
$d.Content.Find.Execute("82-51=31", $true, $true, $false, $false, $false, $true, 1, $false, "26-15=11", 2) | Out-Null
$d.Content.Find.Execute("64-14=50", $true, $true, $false, $false, $false, $true, 1, $false, "80+16=96", 2) | Out-Null
$d.Content.Find.Execute("7+17=24", $true, $true, $false, $false, $false, $true, 1, $false, "30+22=52", 2) | Out-Null
$d.Content.Find.Execute("18+66=84", $true, $true, $false, $false, $false, $true, 1, $false, "18+26=44", 2) | Out-Null
$d.Content.Find.Execute("76-34=42", $true, $true, $false, $false, $false, $true, 1, $false, "88+2=90", 2) | Out-Null
$d.Content.Find.Execute("74-59=15", $true, $true, $false, $false, $false, $true, 1, $false, "55-23=32", 2) | Out-Null
$d.Content.Find.Execute("86+9=95", $true, $true, $false, $false, $false, $true, 1, $false, "52-47=5", 2) | Out-Null
$d.Content.Find.Execute("48+21=69", $true, $true, $false, $false, $false, $true, 1, $false, "44-27=17", 2) | Out-Null
$d.Content.Find.Execute("9+75=84", $true, $true, $false, $false, $false, $true, 1, $false, "32-9=23", 2) | Out-Null
$d.Content.Find.Execute("57-35=22", $true, $true, $false, $false, $false, $true, 1, $false, "3+34=37", 2) | Out-Null
$d.Content.Find.Execute("19-14=5", $true, $true, $false, $false, $false, $true, 1, $false, "78+1=79", 2) | Out-Null
$d.Content.Find.Execute("42+52=94", $true, $true, $false, $false, $false, $true, 1, $false, "22+51=73", 2) | Out-Null
$d.Content.Find.Execute("94-51=43", $true, $true, $false, $false, $false, $true, 1, $false, "36-23=13", 2) | Out-Null
$d.Content.Find.Execute("35+36=71", $true, $true, $false, $false, $false, $true, 1, $false, "29+29=58", 2) | Out-Null
$d.Content.Find.Execute("11-8=3", $true, $true, $false, $false, $false, $true, 1, $false, "10+39=49", 2) | Out-Null
$d.Content.Find.Execute("89-80=9", $true, $true, $false, $false, $false, $true, 1, $false, "60-50=10", 2) | Out-Null
$d.Content.Find.Execute("76-13=63", $true, $true, $false, $false, $false, $true, 1, $false, "48-13=35", 2) | Out-Null
$d.Content.Find.Execute("42+35=77", $true, $true, $false, $false, $false, $true, 1, $false, "31+12=43", 2) | Out-Null
$d.Content.Find.Execute("77-58=19", $true, $true, $false, $false, $false, $true, 1, $false, "50-22=28", 2) | Out-Null
$d.Content.Find.Execute("84-60=24", $true, $true, $false, $false, $false, $true, 1, $false, "57-28=29", 2) | Out-Null
$d.Content.Find.Execute("71-61=10", $true, $true, $false, $false, $false, $true, 1, $false, "5+67=72", 2) | Out-Null
$d.Content.Find.Execute("23-2=21", $true, $true, $false, $false, $false, $true, 1, $false, "72-60=12", 2) | Out-Null
$d.Content.Find.Execute("25+73=98", $true, $true, $false, $false, $false, $true, 1, $false, "16-5=11", 2) | Out-Null
$d.Content.Find.Execute("50-11=39", $true, $true, $false, $false, $false, $true, 1, $false, "10+72=82", 2) | Out-Null
$d.Content.Find.Execute("21+41=62", $true, $true, $false, $false, $false, $true, 1, $false, "22+5=27", 2) | Out-Null
$d.Content.Find.Execute("24+26=50", $true, $true, $false, $false, $false, $true, 1, $false, "24+39=63", 2) | Out-Null
$d.Content.Find.Execute("21+38=59", $true, $true, $false, $false, $false, $true, 1, $false, "15+25=40", 2) | Out-Null
$d.Content.Find.Execute("77-25=52", $true, $true, $false, $false, $false, $true, 1, $false, "90-84=6", 2) | Out-Null
$d.Content.Find.Execute("54+6=60", $true, $true, $false, $false, $false, $true, 1, $false, "8+40=48", 2) | Out-Null
$d.Content.Find.Execute("3+53=56", $true, $true, $false, $false, $false, $true, 1, $false, "13+66=79", 2) | Out-Null
$d.Content.Find.Execute("89-55=34", $true, $true, $false, $false, $false, $true, 1, $false, "26+11=37", 2) | Out-Null
$d.Content.Find.Execute("94-57=37", $true, $true, $false, $false, $false, $true, 1, $false, "21+70=91", 2) | Out-Null
$d.Content.Find.Execute("85-33=52", $true, $true, $false, $false, $false, $true, 1, $false, "4+26=30", 2) | Out-Null
$d.Content.Find.Execute("89-86=3", $true, $true, $false, $false, $false, $true, 1, $false, "82+9=91", 2) | Out-Null
$d.Content.Find.Execute("39-30=9", $true, $true, $false, $false, $false, $true, 1, $false, "27-4=23", 2) | Out-Null
$d.Content.Find.Execute("6+54=60", $true, $true, $false, $false, $false, $true, 1, $false, "84+7=91", 2) | Out-Null
$d.Content.Find.Execute("78-77=1", $true, $true, $false, $false, $false, $true, 1, $false, "74-17=57", 2) | Out-Null
$d.Content.Find.Execute("71-30=41", $true, $true, $false, $false, $false, $true, 1, $false, "4+71=75", 2) | Out-Null
$d.Content.Find.Execute("82-73=9", $true, $true, $false, $false, $false, $true, 1, $false, "12+19=31", 2) | Out-Null
$d.Content.Find.Execute("75-41=34", $true, $true, $false, $false, $false, $true, 1, $false, "38+45=83", 2) | Out-Null
$d.Content.Find.Execute("76-8=68", $true, $true, $false, $false, $false, $true, 1, $false, "56-24=32", 2) | Out-Null
$d.Content.Find.Execute("63+21=84", $true, $true, $false, $false, $false, $true, 1, $false, "33+0=33", 2) | Out-Null
$d.Content.Find.Execute("2+28=30", $true, $true, $false, $false, $false, $true, 1, $false, "24+74=98", 2) | Out-Null
$d.Content.Find.Execute("21-14=7", $true, $true, $false, $false, $false, $true, 1, $false, "44-39=5", 2) | Out-Null
$d.Content.Find.Execute("31+23=54", $true, $true, $false, $false, $false, $true, 1, $false, "80-68=12", 2) | Out-Null
$d.Content.Find.Execute("5+92=97", $true, $true, $false, $false, $false, $true, 1, $false, "79-19=60", 2) | Out-Null
$d.Content.Find.Execute("57-4=53", $true, $true, $false, $false, $false, $true, 1, $false, "54+17=71", 2) | Out-Null
$d.Content.Find.Execute("75-74=1", $true, $true, $false, $false, $false, $true, 1, $false, "23-5=18", 2) | Out-Null
$d.Content.Find.Execute("31+0=31", $true, $true, $false, $false, $false, $true, 1, $false, "5+38=43", 2) | Out-Null
$d.Content.Find.Execute("25+18=43", $true, $true, $false, $false, $false, $true, 1, $false, "35-5=30", 2) | Out-Null
$d.Content.Find.Execute("45+26=71", $true, $true, $false, $false, $false, $true, 1, $false, "85+4=89", 2) | Out-Null
$d.Content.Find.Execute("4+75=79", $true, $true, $false, $false, $false, $true, 1, $false, "27+30=57", 2) | Out-Null
$d.Content.Find.Execute("71-4=67", $true, $true, $false, $false, $false, $true, 1, $false, "71-33=38", 2) | Out-Null
$d.Content.Find.Execute("81+15=96", $true, $true, $false, $false, $false, $true, 1, $false, "13+62=75", 2) | Out-Null
$d.Content.Find.Execute("23+52=75", $true, $true, $false, $false, $false, $true, 1, $false, "43+7=50", 2) | Out-Null
$d.Content.Find.Execute("1+94=95", $true, $true, $false, $false, $false, $true, 1, $false, "61-48=13", 2) | Out-Null
$d.Content.Find.Execute("35+46=81", $true, $true, $false, $false, $false, $true, 1, $false, "50+7=57", 2) | Out-Null
$d.Content.Find.Execute("15+66=81", $true, $true, $false, $false, $false, $true, 1, $false, "38-23=15", 2) | Out-Null
$d.Content.Find.Execute("0+93=93", $true, $true, $false, $false, $false, $true, 1, $false, "39+50=89", 2) | Out-Null
$d.Content.Find.Execute("82-61=21", $true, $true, $false, $false, $false, $true, 1, $false, "36+55=91", 2) | Out-Null
$d.Content.Find.Execute("49+48=97", $true, $true, $false, $false, $false, $true, 1, $false, "44+31=75", 2) | Out-Null
$d.Content.Find.Execute("23-6=17", $true, $true, $false, $false, $false, $true, 1, $false, "21+49=70", 2) | Out-Null
$d.Content.Find.Execute("63-63=0", $true, $true, $false, $false, $false, $true, 1, $false, "83-8=75", 2) | Out-Null
$d.Content.Find.Execute("11+62=73", $true, $true, $false, $false, $false, $true, 1, $false, "19+6=25", 2) | Out-Null
$d.Content.Find.Execute("91-62=29", $true, $true, $false, $false, $false, $true, 1, $false, "76-43=33", 2) | Out-Null
$d.Content.Find.Execute("10+30=40", $true, $true, $false, $false, $false, $true, 1, $false, "69-6=63", 2) | Out-Null
$d.Content.Find.Execute("97-67=30", $true, $true, $false, $false, $false, $true, 1, $false, "21+23=44", 2) | Out-Null
$d.Content.Find.Execute("47+27=74", $true, $true, $false, $false, $false, $true, 1, $false, "92-38=54", 2) | Out-Null
$d.Content.Find.Execute("20+76=96", $true, $true, $false, $false, $false, $true, 1, $false, "61+38=99", 2) | Out-Null
$d.Content.Find.Execute("3+62=65", $true, $true, $false, $false, $false, $true, 1, $false, "46+41=87", 2) | Out-Null
$d.Content.Find.Execute("57+41=98", $true, $true, $false, $false, $false, $true, 1, $false, "23+37=60", 2) | Out-Null
$d.Content.Find.Execute("4+16=20", $true, $true, $false, $false, $false, $true, 1, $false, "56+25=81", 2) | Out-Null
$d.Content.Find.Execute("64-10=54", $true, $true, $false, $false, $false, $true, 1, $false, "63-30=33", 2) | Out-Null
$d.Content.Find.Execute("4+62=66", $true, $true, $false, $false, $false, $true, 1, $false, "58-34=24", 2) | Out-Null
$d.Content.Find.Execute("3+1=4", $true, $true, $false, $false, $false, $true, 1, $false, "87-49=38", 2) | Out-Null
$d.Content.Find.Execute("91-82=9", $true, $true, $false, $false, $false, $true, 1, $false, "87-52=35", 2) | Out-Null
$d.Content.Find.Execute("20+36=56", $true, $true, $false, $false, $false, $true, 1, $false, "70-4=66", 2) | Out-Null
$d.Content.Find.Execute("67-25=42", $true, $true, $false, $false, $false, $true, 1, $false, "22-21=1", 2) | Out-Null
$d.Content.Find.Execute("2+85=87", $true, $true, $false, $false, $false, $true, 1, $false, "13+76=89", 2) | Out-Null
$d.Content.Find.Execute("71+23=94", $true, $true, $false, $false, $false, $true, 1, $false, "7+89=96", 2) | Out-Null
$d.Content.Find.Execute("32+9=41", $true, $true, $false, $false, $false, $true, 1, $false, "69+27=96", 2) | Out-Null
$d.Content.Find.Execute("23+21=44", $true, $true, $false, $false, $false, $true, 1, $false, "65+4=69", 2) | Out-Null
$d.Content.Find.Execute("58-44=14", $true, $true, $false, $false, $false, $true, 1, $false, "95-65=30", 2) | Out-Null
$d.Content.Find.Execute("73-21=52", $true, $true, $false, $false, $false, $true, 1, $false, "76-20=56", 2) | Out-Null
